$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source report stores every Price/Volume cell as literal text (not a
# number). A handful of the new Price values read as plain decimals (e.g.
# '1.00', '0.130'), so without forcing a Text number format first, Excel's
# normal typed-value coercion would store them as numbers and silently drop
# the significant trailing/leading zeros (e.g. '1.00' -> 1, '0.130' -> 0.13).

$ws.Range('D2').Value = '88.245.91'
$ws.Range('E2').Value = '  -1.65%  '
$ws.Range('D3').Value = '3.066.44'
$ws.Range('E3').Value = '  -3.84%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '210.33'
$ws.Range('E5').Value = '  -3.03%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '618.54'
$ws.Range('E6').Value = '  -2.90%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.369'
$ws.Range('E7').Value = '  -5.60%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.801'
$ws.Range('E8').Value = '  +16.11%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('D10').Value = '3.063.59'
$ws.Range('E10').Value = '  -3.71%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.592'
$ws.Range('E11').Value = '  +3.57%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.180'
$ws.Range('E12').Value = '  +0.15%  '
$ws.Range('E13').Value = '  -6.96%  '
$ws.Range('E14').Value = '  -1.74%  '
$ws.Range('D15').Value = '87.974.31'
$ws.Range('E15').Value = '  -1.70%  '
$ws.Range('B16').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C16').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D16').Value = '3.633.45'
$ws.Range('E16').Value = '  -3.70%  '
$ws.Range('B17').Value = 'Avalanche'
$ws.Range('C17').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '31.85'
$ws.Range('E17').Value = '  -4.10%  '
$ws.Range('D18').Value = '3.057.69'
$ws.Range('E18').Value = '  -4.90%  '
$ws.Range('E19').Value = '  -3.79%  '
$ws.Range('E20').Value = '  -10.64%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.26'
$ws.Range('E21').Value = '  -1.63%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '420.47'
$ws.Range('E22').Value = '  -3.56%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '8.15'
$ws.Range('E23').Value = '  -5.72%  '
$ws.Range('E24').Value = '  -4.29%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '5.45'
$ws.Range('E25').Value = '  +2.76%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '11.76'
$ws.Range('E26').Value = '  -1.47%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '82.08'
$ws.Range('E27').Value = '  +0.70%  '
$ws.Range('D28').Value = '3.231.99'
$ws.Range('E28').Value = '  -3.91%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('E30').Value = '  +8.61%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.170'
$ws.Range('E31').Value = '  +6.72%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '8.01'
$ws.Range('E32').Value = '  -5.25%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '508.04'
$ws.Range('E33').Value = '  -6.94%  '
$ws.Range('E34').Value = '  -11.25%  '
$ws.Range('E35').Value = '  -4.31%  '
$ws.Range('E36').Value = '  -5.75%  '
$ws.Range('E37').Value = '  -6.65%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '22.24'
$ws.Range('E38').Value = '  -0.70%  '
$ws.Range('B39').Value = 'WhiteBITCoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '22.21'
$ws.Range('E39').Value = '  -0.77%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.130'
$ws.Range('E40').Value = '  +1.05%  '
$ws.Range('E41').Value = '  +0.22%  '
$ws.Range('E43').Value = '  -3.72%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '148.46'
$ws.Range('E44').Value = '  +1.63%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.80'
$ws.Range('E45').Value = '  -6.55%  '
$ws.Range('E46').Value = '  +6.05%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '43.36'
$ws.Range('E47').Value = '  -0.68%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0688'
$ws.Range('E48').Value = '  +13.91%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '157.22'
$ws.Range('E49').Value = '  -9.44%  '
$ws.Range('B50').Value = 'ImmutableX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.18'
$ws.Range('E50').Value = '  -4.95%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.701'
$ws.Range('E51').Value = '  -7.42%  '
